$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 45192.45601851852
$ws.Cells.Item(3, 3).Value = 45192.45664351852
$ws.Cells.Item(3, 4).Value = 'Hossam.Ibrahim'
$ws.Cells.Item(3, 5).Value = 'Hossam Tabana'
$ws.Cells.Item(3, 6).Value = ""
$ws.Cells.Item(3, 7).Value = ""
$ws.Cells.Item(3, 8).Value = ""
$ws.Cells.Item(3, 9).Value = 'Geospatial Maps'
$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(3, 11).Value = ""
$ws.Cells.Item(3, 12).Value = 'Sometimes'
$ws.Cells.Item(3, 13).Value = ""
$ws.Cells.Item(3, 14).Value = ""
$ws.Cells.Item(3, 15).Value = 'Yes, simple calculations'
$ws.Cells.Item(3, 16).Value = ""
$ws.Cells.Item(3, 17).Value = ""
$ws.Cells.Item(3, 18).Value = 'Beginner level'
$ws.Cells.Item(3, 19).Value = ""
$ws.Cells.Item(3, 20).Value = ""
$ws.Cells.Item(3, 21).Value = 'Basic transformations only'
$ws.Cells.Item(3, 22).Value = ""
$ws.Cells.Item(3, 23).Value = ""
$ws.Cells.Item(3, 24).Value = 'Yes, basic automation'
$ws.Cells.Item(3, 25).Value = ""
$ws.Cells.Item(3, 26).Value = ""
$ws.Cells.Item(3, 27).Value = 'Occasionally'
$ws.Cells.Item(3, 28).Value = ""
$ws.Cells.Item(3, 29).Value = ""
$ws.Cells.Item(3, 30).Value = 'Very Important'
$ws.Cells.Item(3, 31).Value = ""
$ws.Cells.Item(3, 32).Value = ""
$ws.Cells.Item(3, 33).Value = 'Yes, as static files'
$ws.Cells.Item(3, 34).Value = ""
$ws.Cells.Item(3, 35).Value = ""
$ws.Cells.Item(3, 36).Value = 'Occasionally'
$ws.Cells.Item(3, 37).Value = ""
$ws.Cells.Item(3, 38).Value = ""
$ws.Cells.Item(3, 39).Value = 'Yes, it''s essential'
$ws.Cells.Item(3, 40).Value = ""
$ws.Cells.Item(3, 41).Value = ""
$ws.Cells.Item(3, 42).Value = 'Not concerned'
$ws.Cells.Item(3, 43).Value = ""
$ws.Cells.Item(3, 44).Value = ""
$ws.Cells.Item(3, 45).Value = 'Just exploring'
$ws.Cells.Item(3, 46).Value = ""
$ws.Cells.Item(3, 47).Value = ""
$ws.Cells.Item(3, 48).Value = 'Yes, basic trend lines'
$ws.Cells.Item(3, 49).Value = ""
$ws.Cells.Item(3, 50).Value = ""
$ws.Cells.Item(3, 51).Value = 'Beginner'
$ws.Cells.Item(3, 52).Value = ""
$ws.Cells.Item(3, 53).Value = ""
$ws.Cells.Item(3, 54).Value = 'Definitely'
$ws.Cells.Item(3, 55).Value = ""
$ws.Cells.Item(3, 56).Value = ""
$ws.Cells.Item(3, 57).Value = 'I know what it is but haven''t used it'
$ws.Cells.Item(3, 58).Value = ""
$ws.Cells.Item(3, 59).Value = ""
$ws.Cells.Item(3, 60).Value = 'Possibly'
$ws.Cells.Item(3, 61).Value = ""
$ws.Cells.Item(3, 62).Value = ""
$ws.Cells.Item(3, 63).Value = 'Yes, to multiple formats'
$ws.Cells.Item(3, 64).Value = ""
$ws.Cells.Item(3, 65).Value = ""
$ws.Cells.Item(3, 66).Value = 'Using Power BI workspaces'
$ws.Cells.Item(3, 67).Value = ""
$ws.Cells.Item(3, 68).Value = ""
$ws.Cells.Item(3, 69).Value = 'Advanced Training'

$ws.Cells.Item(3, 2).NumberFormat = $ws.Cells.Item(2, 2).NumberFormat
$ws.Cells.Item(3, 3).NumberFormat = $ws.Cells.Item(2, 3).NumberFormat
